$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 51005.00823254297
$ws.Range("B3").Value = 32401.14538257632
$ws.Range("B5").Value = 2081.291180778615
$ws.Range("B6").Value = 3532.226790597309
$ws.Range("B8").Value = 1964.726790597309
$ws.Range("B9").Value = -13.85159912109373
$ws.Range("B10").Value = 26.99999999998431
$ws.Range("B11").Value = 26.62374482422558
$ws.Range("B12").Value = 27.03821225603286
$ws.Range("B13").Value = 64.73487797803335
$ws.Range("B14").Value = 11.63509226463992
$ws.Range("B15").Value = 8.422298003455047
$ws.Range("B16").Value = 15.45561770905346
$ws.Range("B17").Value = 15.05671624322933
$ws.Range("B18").Value = 1.214054630626317
$ws.Range("B19").Value = 1.99999999999919
$ws.Range("B20").Value = 0.218895274988578
$ws.Range("B21").Value = 1.816104725010612
$ws.Range("B22").Value = 1.819294851080049
$ws.Range("B23").Value = 0.2189538329839706
$ws.Range("B25").Value = 1.038248658180237
$ws.Range("B26").Value = 0.1228985987824348
$ws.Range("B29").Value = 104.5975625
$ws.Range("B31").Value = 104.5975605534461
